# Update the non-functional requirement descriptions (NR1 and NR2 text),
# as captured by the commit "funkcni & nefunkcni pozadavky".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hárok1")

# NR1 description: "upgradeable" -> "scalable"
$ws.Range("I4").Value = "Aplikace bude dobře škálovatelná"

# NR2 description: "handle influx of users" -> "handle heavy load of users"
$ws.Range("I5").Value = "Aplikace bude zvládat velkou zátěž uživatelů"

# NR3 text stays the same ("Aplikace bude mít podporu na příštích 5 let.")

# Reflect the last active cell selection at save time
$ws.Range("C8").Select()
